$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.064.96'
$ws.Range("E2").Value = '  +3.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.254.56'
$ws.Range("E3").Value = '  +2.22%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.26'
$ws.Range("E5").Value = '  +2.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '79.45'
$ws.Range("E6").Value = '  +7.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +2.24%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  +3.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.24'
$ws.Range("E10").Value = '  +7.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0930'
$ws.Range("E11").Value = '  +1.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.11'
$ws.Range("E12").Value = '  +4.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  +1.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.591.21'
$ws.Range("E14").Value = '  +1.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.75'
$ws.Range("E15").Value = '  +2.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.262.98'
$ws.Range("E16").Value = '  +1.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.796'
$ws.Range("E17").Value = '  +2.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.964.69'
$ws.Range("E18").Value = '  +3.19%  '

$ws.Range("E19").Value = '  +1.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.63'
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.08'
$ws.Range("E21").Value = '  +3.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.36'
$ws.Range("E22").Value = '  +8.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.19'
$ws.Range("E23").Value = '  +2.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.60'
$ws.Range("E24").Value = '  +2.19%  '

$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '42.37'
$ws.Range("E26").Value = '  +8.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.86'
$ws.Range("E27").Value = '  +1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.38'
$ws.Range("E28").Value = '  -0.15%  '

$ws.Range("E29").Value = '  +1.78%  '

$ws.Range("E30").Value = '  -0.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.68'
$ws.Range("E31").Value = '  +2.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.71'
$ws.Range("E32").Value = '  +3.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0878'
$ws.Range("E33").Value = '  +10.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.32'
$ws.Range("E34").Value = '  +2.56%  '

$ws.Range("E35").Value = '  +2.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("E36").Value = '  +4.86%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0362'
$ws.Range("E37").Value = '  +11.39%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.51'
$ws.Range("E38").Value = '  +3.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.56'
$ws.Range("E39").Value = '  +13.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.89'
$ws.Range("E40").Value = '  +20.11%  '

$ws.Range("E41").Value = '  +3.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '62.71'
$ws.Range("E42").Value = '  +6.79%  '

$ws.Range("E43").Value = '  +2.64%  '

$ws.Range("E44").Value = '  +2.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.33'
$ws.Range("E45").Value = '  +2.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.58'
$ws.Range("E46").Value = '  +0.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.473'
$ws.Range("E47").Value = '  -0.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0991'
$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("E49").Value = '  +2.34%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.52'
$ws.Range("E50").Value = '  +26.97%  '

$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
$ws.Range("E51").Value = '  +2.02%  '
